$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = "2024 - Vår"
$ws.Range("B16").Value = "[Oppgaveformulering](tidligere-eksamensoppgaver/skole-24-v.pdf)"
$ws.Range("C16").Value = "[Løsningsforslag](tidligere-eksamensoppgaver/skole-24-v-fasit.pdf)"

$ws.Range("B17").Select()
